$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.399.56'
$ws.Cells.Item(2, 5).Value = '  +1.07%  '
$ws.Cells.Item(3, 4).Value = '2.376.23'
$ws.Cells.Item(3, 5).Value = '  +3.44%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).Value = "'309.84"
$ws.Cells.Item(5, 5).Value = '  +0.07%  '
$ws.Cells.Item(6, 4).Value = "'105.21"
$ws.Cells.Item(6, 5).Value = '  +5.32%  '
$ws.Cells.Item(7, 4).Value = "'0.523"
$ws.Cells.Item(7, 5).Value = '  -2.20%  '
$ws.Cells.Item(8, 5).Value = '  -0.01%  '
$ws.Cells.Item(9, 4).Value = "'0.519"
$ws.Cells.Item(9, 5).Value = '  +0.48%  '
$ws.Cells.Item(10, 4).Value = "'36.22"
$ws.Cells.Item(10, 5).Value = '  +0.70%  '
$ws.Cells.Item(11, 4).Value = "'53.44"
$ws.Cells.Item(11, 5).Value = '  +2.62%  '
$ws.Cells.Item(12, 5).Value = '  -1.00%  '
$ws.Cells.Item(13, 4).Value = "'0.112"
$ws.Cells.Item(13, 5).Value = '  -0.51%  '
$ws.Cells.Item(14, 4).Value = "'7.02"
$ws.Cells.Item(14, 5).Value = '  -1.16%  '
$ws.Cells.Item(15, 4).Value = '2.744.14'
$ws.Cells.Item(15, 5).Value = '  +3.42%  '
$ws.Cells.Item(16, 4).Value = "'15.65"
$ws.Cells.Item(16, 5).Value = '  +5.44%  '
$ws.Cells.Item(17, 4).Value = '2.370.94'
$ws.Cells.Item(17, 5).Value = '  +2.59%  '
$ws.Cells.Item(18, 5).Value = '  +2.13%  '
$ws.Cells.Item(19, 4).Value = '43.364.33'
$ws.Cells.Item(19, 5).Value = '  +1.18%  '
$ws.Cells.Item(20, 5).Value = '  -3.54%  '
$ws.Cells.Item(21, 5).Value = '  -0.14%  '
$ws.Cells.Item(22, 4).Value = "'6.27"
$ws.Cells.Item(22, 5).Value = '  +3.46%  '
$ws.Cells.Item(23, 5).Value = '  +0.48%  '
$ws.Cells.Item(24, 4).Value = "'242.00"
$ws.Cells.Item(24, 5).Value = '  +1.13%  '
$ws.Cells.Item(25, 5).Value = '  +2.38%  '
$ws.Cells.Item(26, 4).Value = "'2.63"
$ws.Cells.Item(26, 5).Value = '  +0.64%  '
$ws.Cells.Item(27, 5).Value = '  +0.14%  '
$ws.Cells.Item(28, 4).Value = "'25.94"
$ws.Cells.Item(28, 5).Value = '  +7.77%  '
$ws.Cells.Item(29, 2).Value = 'LEO'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(29, 4).Value = "'3.81"
$ws.Cells.Item(29, 5).Value = '  -4.77%  '
$ws.Cells.Item(30, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(30, 4).Value = "'37.15"
$ws.Cells.Item(30, 5).Value = '  -3.53%  '
$ws.Cells.Item(31, 2).Value = 'Cosmos'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(31, 4).Value = "'9.61"
$ws.Cells.Item(31, 5).Value = '  -0.14%  '
$ws.Cells.Item(32, 2).Value = 'Toncoin'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(32, 4).Value = "'2.12"
$ws.Cells.Item(32, 5).Value = '  +0.05%  '
$ws.Cells.Item(33, 2).Value = 'Monero'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(33, 4).Value = "'161.88"
$ws.Cells.Item(33, 5).Value = '  -3.62%  '
$ws.Cells.Item(34, 2).Value = 'Filecoin'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(34, 4).Value = "'5.28"
$ws.Cells.Item(34, 5).Value = '  -0.72%  '
$ws.Cells.Item(35, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(35, 4).Value = "'1.00"
$ws.Cells.Item(35, 5).Value = '  +0.04%  '
$ws.Cells.Item(36, 2).Value = 'Celestia'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Cells.Item(36, 4).Value = "'18.28"
$ws.Cells.Item(36, 5).Value = '  +3.19%  '
$ws.Cells.Item(37, 4).Value = "'4.79"
$ws.Cells.Item(37, 5).Value = '  +13.02%  '
$ws.Cells.Item(38, 2).Value = 'WEMIXToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(38, 4).Value = "'2.54"
$ws.Cells.Item(38, 5).Value = '  +6.54%  '
$ws.Cells.Item(39, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(39, 4).Value = "'3.13"
$ws.Cells.Item(39, 5).Value = '  +0.42%  '
$ws.Cells.Item(40, 2).Value = 'Hedera'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(40, 4).Value = "'0.0744"
$ws.Cells.Item(40, 5).Value = '  +0.97%  '
$ws.Cells.Item(41, 2).Value = 'ARBITRUM'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(41, 4).Value = "'1.95"
$ws.Cells.Item(41, 5).Value = '  +6.40%  '
$ws.Cells.Item(42, 2).Value = 'Kaspa'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(42, 4).Value = "'0.106"
$ws.Cells.Item(42, 5).Value = '  +0.97%  '
$ws.Cells.Item(43, 2).Value = 'Stellar'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(43, 4).Value = "'0.114"
$ws.Cells.Item(43, 5).Value = '  -1.30%  '
$ws.Cells.Item(44, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(44, 4).Value = "'2.48"
$ws.Cells.Item(44, 5).Value = '  +9.01%  '
$ws.Cells.Item(45, 2).Value = 'EnergySwap'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(45, 4).Value = "'19.94"
$ws.Cells.Item(45, 5).Value = '  +3.66%  '
$ws.Cells.Item(46, 2).Value = 'Maker'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(46, 4).Value = '2.006.14'
$ws.Cells.Item(46, 5).Value = '  +2.14%  '
$ws.Cells.Item(47, 2).Value = 'VeChain'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(47, 4).Value = "'0.0290"
$ws.Cells.Item(47, 5).Value = '  +0.93%  '
$ws.Cells.Item(48, 2).Value = 'NEARProtocol'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(48, 4).Value = "'3.16"
$ws.Cells.Item(48, 5).Value = '  +5.47%  '
$ws.Cells.Item(49, 2).Value = 'FraxShare'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(49, 4).Value = "'10.53"
$ws.Cells.Item(49, 5).Value = '  +7.65%  '
$ws.Cells.Item(50, 2).Value = 'MultiversX'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Cells.Item(50, 4).Value = "'58.25"
$ws.Cells.Item(50, 5).Value = '  +5.92%  '
$ws.Cells.Item(51, 2).Value = 'HuobiToken'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(51, 4).Value = "'2.96"
$ws.Cells.Item(51, 5).Value = '  +0.81%  '
